# Updates the cryptos price/volume table with refreshed values.
# Numeric-looking Price (column D) values are written with a leading
# apostrophe so Excel keeps them as text (matching the source data,
# which stores prices like "431.70" / "0.999" as literal strings
# rather than numbers that would lose trailing zeros).
$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range('D2').Value = '71.892.03'
$ws.Range('E2').Value = '  -0.66%  '
$ws.Range('D3').Value = '3.995.18'
$ws.Range('E3').Value = '  -1.05%  '
$ws.Range('D4').Value = '''0.999'
$ws.Range('E4').Value = '  -0.20%  '
$ws.Range('D5').Value = '''544.61'
$ws.Range('E5').Value = '  +4.68%  '
$ws.Range('D6').Value = '''150.81'
$ws.Range('E6').Value = '  +1.92%  '
$ws.Range('E7').Value = '  +10.70%  '
$ws.Range('E8').Value = '  +0.08%  '
$ws.Range('E9').Value = '  +0.31%  '
$ws.Range('E10').Value = '  -3.28%  '
$ws.Range('D11').Value = '''53.29'
$ws.Range('E11').Value = '  +12.49%  '
$ws.Range('D12').Value = '''0.0000324'
$ws.Range('E12').Value = '  -3.24%  '
$ws.Range('D13').Value = '''10.68'
$ws.Range('E13').Value = '  -2.16%  '
$ws.Range('D14').Value = '4.637.51'
$ws.Range('E14').Value = '  -1.01%  '
$ws.Range('D15').Value = '4.005.21'
$ws.Range('E15').Value = '  -0.64%  '
$ws.Range('D16').Value = '''14.13'
$ws.Range('E16').Value = '  -0.96%  '
$ws.Range('E17').Value = '  -3.40%  '
$ws.Range('E18').Value = '  -0.22%  '
$ws.Range('E19').Value = '  -2.02%  '
$ws.Range('D20').Value = '71.869.58'
$ws.Range('E20').Value = '  -0.66%  '
$ws.Range('D21').Value = '''431.70'
$ws.Range('E21').Value = '  -1.78%  '
$ws.Range('D22').Value = '''96.85'
$ws.Range('E22').Value = '  -4.38%  '
$ws.Range('E23').Value = '  -0.51%  '
$ws.Range('D24').Value = '''4.29'
$ws.Range('E24').Value = '  +6.41%  '
$ws.Range('D25').Value = '''14.34'
$ws.Range('E25').Value = '  -2.77%  '
$ws.Range('E26').Value = '  -2.43%  '
$ws.Range('D27').Value = '''10.72'
$ws.Range('E27').Value = '  -4.97%  '
$ws.Range('D28').Value = '''5.85'
$ws.Range('E28').Value = '  +1.03%  '
$ws.Range('D29').Value = '''36.75'
$ws.Range('E29').Value = '  -2.10%  '
$ws.Range('D30').Value = '''3.63'
$ws.Range('E30').Value = '  +17.41%  '
$ws.Range('D31').Value = '''7.50'
$ws.Range('E31').Value = '  +8.16%  '
$ws.Range('B32').Value = 'Cosmos'
$ws.Range('C32').Value = 'https://coinranking.com/coin/Knsels4_Ol-Ny+cosmos-atom'
$ws.Range('D32').Value = '''13.45'
$ws.Range('E32').Value = '  -0.85%  '
$ws.Range('B33').Value = 'Hedera'
$ws.Range('C33').Value = 'https://coinranking.com/coin/jad286TjB+hedera-hbar'
$ws.Range('D33').Value = '''0.131'
$ws.Range('E33').Value = '  +1.57%  '
$ws.Range('D34').Value = '''48.96'
$ws.Range('E34').Value = '  +17.57%  '
$ws.Range('D35').Value = '''679.14'
$ws.Range('E35').Value = '  -2.03%  '
$ws.Range('D36').Value = '''65.86'
$ws.Range('E36').Value = '  -3.30%  '
$ws.Range('E37').Value = '  +0.35%  '
$ws.Range('B38').Value = 'PEPE'
$ws.Range('C38').Value = 'https://coinranking.com/coin/03WI8NQPF+pepe-pepe'
$ws.Range('D38').Value = '0.0₃0831'
$ws.Range('E38').Value = '  -6.45%  '
$ws.Range('B39').Value = 'Kaspa'
$ws.Range('C39').Value = 'https://coinranking.com/coin/V8GxkwWow+kaspa-kas'
$ws.Range('D39').Value = '''0.152'
$ws.Range('E39').Value = '  -1.00%  '
$ws.Range('D40').Value = '''3.39'
$ws.Range('E40').Value = '  -6.95%  '
$ws.Range('D41').Value = '''3.35'
$ws.Range('E41').Value = '  +5.14%  '
$ws.Range('E42').Value = '  +0.09%  '
$ws.Range('E43').Value = '  +0.25%  '
$ws.Range('D44').Value = '''0.0487'
$ws.Range('E44').Value = '  -0.75%  '
$ws.Range('D45').Value = '''2.78'
$ws.Range('E45').Value = '  -0.92%  '
$ws.Range('E46').Value = '  +1.42%  '
$ws.Range('D47').Value = '''9.82'
$ws.Range('E47').Value = '  +8.08%  '
$ws.Range('D48').Value = '''3.38'
$ws.Range('E48').Value = '  -3.53%  '
$ws.Range('D49').Value = '''0.000280'
$ws.Range('E49').Value = '  +1.91%  '
$ws.Range('D50').Value = '''2.99'
$ws.Range('E50').Value = '  -4.07%  '
$ws.Range('D51').Value = '''144.76'
$ws.Range('E51').Value = '  +1.37%  '
